$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Description")

# --- Text updates as supplied by PM&C ---

# Influences: drop the inline citation; the source now lives in the new
# References block below instead of being parenthesised inline.
$ws.Range("B7").Value = "Tobacco smoking is a multi-determined behaviour, influenced by a range of biological, psychosocial, and environmental factors."
$ws.Rows.Item(7).RowHeight = 25.35

# Notes: replace the old one-line source note with the NT data-caveat note.
$ws.Range("B8").Value = "Data for the Northern Territory should be used with care as very remote areas were excluded from the Australian Health Survey, which translates to exclusion of around 25 per cent of the Northern Territory population. " + [char]10
# The note now reads as body copy, so give it the same look as the rest of
# the description column (B3:B7) instead of the old plain footnote style.
$ws.Range("B8").Style = $ws.Range("B7").Style
$ws.Rows.Item(8).RowHeight = 49.45

# --- New "References" block (rows 9-11) ---
$ws.Range("A9").Value = "References"
$ws.Range("B9").Value = "ABS (unpublished) National Health Survey 2014-15"
$ws.Range("B10").Value = "ABS (unpublished) Australian Health Survey 2011" + [char]173 + "13 (2011" + [char]173 + "12 core component)"
$ws.Range("B11").Value = "ABS (unpublished) National Health Survey 2007-08."

# The reference lines pick up the small plain-Arial footnote look the old
# Notes cell (B8) used to carry.
$ws.Range("B9:B11").Font.Name = "Arial"
$ws.Range("B9:B11").Font.Size = 10
$ws.Range("B9:B11").WrapText = $true

$ws.Rows.Item(9).RowHeight = 12.8
$ws.Rows.Item(10).RowHeight = 12.8
$ws.Rows.Item(11).RowHeight = 12.8
